# Commit: "commit a basic project framework"
#
# The "words" worksheet had a duplicate/incomplete entry for the word
# "assemble" (phonetic transcription, Chinese definition, and example
# sentence) that needs to be removed, leaving only the headword in
# column A (row 4) and the "已处理" status in column E. Also update the
# active selection to reflect the cells that were just cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("words")

# Remove the phonetic/definition/example text for the "assemble" row,
# keeping the word itself (A4) and processing status (E4) untouched.
$ws.Range("B4:D4").ClearContents()

# Reflect the edited range as the active selection, as was left after
# the edit.
$ws.Activate()
$ws.Range("B4:D4").Select()
